# Apply fix based on feedback
#
# Change the closing sentence of the "Group" text box on slide 1 from
# "...return with the last failure." to
# "...return with the last result or exception."
#
# The run containing " failure." is the 3rd run of paragraph 7 in the
# shape named "Group" (id 314), but we locate it defensively by scanning
# for the exact run text instead of relying purely on fixed indices.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetOld = " failure."
$targetNew = " result or exception."

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }

    $tr = $sh.TextFrame.TextRange
    $paraCount = $tr.Paragraphs().Count

    for ($pIdx = 1; $pIdx -le $paraCount; $pIdx++) {
        $para = $tr.Paragraphs($pIdx)
        $runCount = $para.Runs().Count

        for ($rIdx = 1; $rIdx -le $runCount; $rIdx++) {
            $run = $para.Runs($rIdx)
            if ($run.Text -eq $targetOld) {
                $run.Text = $targetNew
            }
        }
    }
}
